$d = $word.ActiveDocument

# Locate the first bullet paragraph of the EvolutionIQ job ("Orchestrated Google
# Cloud infrastructure...") and insert a brand-new bullet paragraph immediately
# before it. Anchoring the insertion on the following bullet (rather than the
# preceding "New York City, NY - Feb 2024 - Current" subtitle line) makes the
# new paragraph inherit the correct "Normal" style plus the existing numbered
# list formatting (numId 3, ilvl 0) used by the rest of that bullet list.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Orchestrated*Google Cloud*infrastructure*") {
        $p.Range.InsertParagraphBefore()
        break
    }
}

# The insertion above left a new, empty paragraph directly before the
# "Orchestrated..." bullet. Find it (it now sits right after the "New York
# City, NY - Feb 2024 - Current" line and is still empty) and give it its text.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "`r" -or $t -eq "") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text -like "*New York City, NY - Feb 2024 - Current*") {
            $p.Range.Text = "Created infrastructure to provide ephemeral environments, allowing testing of code branches before sending those code changes to the company at large."
            break
        }
    }
}
